$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every edited cell in the source data is stored as literal text
# (inline string), even when the text looks like a plain number
# (e.g. "317.09") or uses a locale-style separator (e.g. "28.017.20").
# Force text format first so Excel does not silently reinterpret
# the assigned string as a numeric value and reformat/round it.
$editedCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "D5",
    "E5",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "D10",
    "E10",
    "E11",
    "E12",
    "D13",
    "E13",
    "E14",
    "B15",
    "C15",
    "D15",
    "E15",
    "B16",
    "C16",
    "D16",
    "E16",
    "E17",
    "E18",
    "D19",
    "E19",
    "E20",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D28",
    "E28",
    "D29",
    "E29",
    "D30",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "B35",
    "C35",
    "D35",
    "E35",
    "B36",
    "C36",
    "D36",
    "E36",
    "D37",
    "E37",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "B41",
    "C41",
    "D41",
    "E41",
    "B42",
    "C42",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($ec in $editedCells) {
    $ws.Range($ec).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.017.20"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "1.791.76"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "317.09"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.5364"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("D8").Value = "0.3770"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").Value = "0.07423"
$ws.Range("D10").Value = "41.75"
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "20.54"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.789.72"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.232"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "0.06493"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "5.892"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "28.039.25"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").Value = "2.089"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "155.56"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("D27").Value = "20.30"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "1.992.92"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "2.304"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").Value = "121.09"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Value = "1.118"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "0.1060"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("D33").Value = "3.655"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "5.553"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.06502"
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D36").Value = "0.2250"
$ws.Range("E36").Value = "  -5.38%  "
$ws.Range("D37").Value = "0.02289"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").Value = "8.469"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").Value = "0.6181"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "1.447"
$ws.Range("E41").Value = "  +2.95%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "11.12"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("D43").Value = "1.176"
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "13.25"
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("D46").Value = "3.670"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "0.5770"
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("D48").Value = "124.80"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "1.185"
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("D50").Value = "1.923"
$ws.Range("E50").Value = "  -3.93%  "
$ws.Range("D51").Value = "0.06820"
$ws.Range("E51").Value = "  -1.72%  "

# Restore the default (unstyled) cell format so only values changed,
# matching the source workbook where these cells carry no explicit
# number format.
foreach ($ec in $editedCells) {
    $ws.Range($ec).Style = "Normal"
}
